$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.70913233333333
$ws.Range("H2").Value = 32.127397
$ws.Range("I2").Value = 0.007451729107954897
$ws.Range("J2").Value = 0.007451729107954897
$ws.Range("M2").Value = 5.854382333333334
$ws.Range("N2").Value = 17.563147
$ws.Range("O2").Value = 0.1730451459016118
$ws.Range("P2").Value = 0.1730451459016118
$ws.Range("Q2").Value = 62.69535513759546
$ws.Range("R2").Value = 564.258196238359
$ws.Range("S2").Value = 0.001289485550705343
$ws.Range("T2").Value = 0.001289485550705343
$ws.Range("G3").Value = 10.70913233333333
$ws.Range("H3").Value = 32.127397
$ws.Range("I3").Value = 0.007451729107954897
$ws.Range("J3").Value = 0.007451729107954897
$ws.Range("O3").Value = 0.1208497063316524
$ws.Range("P3").Value = 0.1208497063316525
$ws.Range("Q3").Value = 43.78461595822489
$ws.Range("R3").Value = 394.061543624024
$ws.Range("S3").Value = 0.0009005392743593757
$ws.Range("T3").Value = 0.0009005392743593758
$ws.Range("G4").Value = 10.70913233333333
$ws.Range("H4").Value = 32.127397
$ws.Range("I4").Value = 0.007451729107954897
$ws.Range("J4").Value = 0.007451729107954897
$ws.Range("M4").Value = 11.64342866666667
$ws.Range("N4").Value = 34.930286
$ws.Range("O4").Value = 0.3441590756630932
$ws.Range("P4").Value = 0.3441590756630932
$ws.Range("Q4").Value = 124.6910184050602
$ws.Range("R4").Value = 1122.219165645542
$ws.Range("S4").Value = 0.002564580201885523
$ws.Range("T4").Value = 0.002564580201885523
$ws.Range("G5").Value = 10.70913233333333
$ws.Range("H5").Value = 32.127397
$ws.Range("I5").Value = 0.007451729107954897
$ws.Range("J5").Value = 0.007451729107954897
$ws.Range("M5").Value = 0.1645376666666667
$ws.Range("N5").Value = 0.493613
$ws.Range("O5").Value = 0.004863441250245888
$ws.Range("P5").Value = 0.004863441250245888
$ws.Range("Q5").Value = 1.762055646151222
$ws.Range("R5").Value = 15.858500815361
$ws.Range("S5").Value = 0.00003624104672928584
$ws.Range("T5").Value = 0.00003624104672928584
$ws.Range("G6").Value = 10.70913233333333
$ws.Range("H6").Value = 32.127397
$ws.Range("I6").Value = 0.007451729107954897
$ws.Range("J6").Value = 0.007451729107954897
$ws.Range("M6").Value = 12.08065233333333
$ws.Range("N6").Value = 36.241957
$ws.Range("O6").Value = 0.3570826308533967
$ws.Range("P6").Value = 0.3570826308533967
$ws.Range("Q6").Value = 129.3733045106588
$ws.Range("R6").Value = 1164.359740595929
$ws.Range("S6").Value = 0.00266088303427537
$ws.Range("T6").Value = 0.00266088303427537
$ws.Range("I7").Value = 0.03290895798513831
$ws.Range("J7").Value = 0.03290895798513832
$ws.Range("M7").Value = 5.854382333333334
$ws.Range("N7").Value = 17.563147
$ws.Range("O7").Value = 0.1730451459016118
$ws.Range("P7").Value = 0.1730451459016118
$ws.Range("Q7").Value = 276.8805438571162
$ws.Range("R7").Value = 2491.924894714045
$ws.Range("S7").Value = 0.005694735436008273
$ws.Range("T7").Value = 0.005694735436008274
$ws.Range("I8").Value = 0.03290895798513831
$ws.Range("J8").Value = 0.03290895798513832
$ws.Range("O8").Value = 0.1208497063316524
$ws.Range("P8").Value = 0.1208497063316525
$ws.Range("S8").Value = 0.003977037908184654
$ws.Range("T8").Value = 0.003977037908184655
$ws.Range("I9").Value = 0.03290895798513831
$ws.Range("J9").Value = 0.03290895798513832
$ws.Range("M9").Value = 11.64342866666667
$ws.Range("N9").Value = 34.930286
$ws.Range("O9").Value = 0.3441590756630932
$ws.Range("P9").Value = 0.3441590756630932
$ws.Range("Q9").Value = 550.6710491442454
$ws.Range("R9").Value = 4956.039442298209
$ws.Range("S9").Value = 0.01132591656120077
$ws.Range("T9").Value = 0.01132591656120077
$ws.Range("I10").Value = 0.03290895798513831
$ws.Range("J10").Value = 0.03290895798513832
$ws.Range("M10").Value = 0.1645376666666667
$ws.Range("N10").Value = 0.493613
$ws.Range("O10").Value = 0.004863441250245888
$ws.Range("P10").Value = 0.004863441250245888
$ws.Range("Q10").Value = 7.781739564950556
$ws.Range("R10").Value = 70.035656084555
$ws.Range("S10").Value = 0.0001600507837675305
$ws.Range("T10").Value = 0.0001600507837675305
$ws.Range("I11").Value = 0.03290895798513831
$ws.Range("J11").Value = 0.03290895798513832
$ws.Range("M11").Value = 12.08065233333333
$ws.Range("N11").Value = 36.241957
$ws.Range("O11").Value = 0.3570826308533967
$ws.Range("P11").Value = 0.3570826308533967
$ws.Range("Q11").Value = 571.3493580965994
$ws.Range("R11").Value = 5142.144222869395
$ws.Range("S11").Value = 0.01175121729597709
$ws.Range("T11").Value = 0.01175121729597709
$ws.Range("G12").Value = 411.37678
$ws.Range("H12").Value = 1234.13034
$ws.Range("I12").Value = 0.2862480573072345
$ws.Range("J12").Value = 0.2862480573072345
$ws.Range("M12").Value = 5.854382333333334
$ws.Range("N12").Value = 17.563147
$ws.Range("O12").Value = 0.1730451459016118
$ws.Range("P12").Value = 0.1730451459016118
$ws.Range("Q12").Value = 2408.356953175553
$ws.Range("R12").Value = 21675.21257857998
$ws.Range("S12").Value = 0.04953383684078333
$ws.Range("T12").Value = 0.04953383684078333
$ws.Range("G13").Value = 411.37678
$ws.Range("H13").Value = 1234.13034
$ws.Range("I13").Value = 0.2862480573072345
$ws.Range("J13").Value = 0.2862480573072345
$ws.Range("O13").Value = 0.1208497063316524
$ws.Range("P13").Value = 0.1208497063316525
$ws.Range("Q13").Value = 1681.926580584586
$ws.Range("R13").Value = 15137.33922526128
$ws.Range("S13").Value = 0.0345929936635853
$ws.Range("T13").Value = 0.03459299366358531
$ws.Range("G14").Value = 411.37678
$ws.Range("H14").Value = 1234.13034
$ws.Range("I14").Value = 0.2862480573072345
$ws.Range("J14").Value = 0.2862480573072345
$ws.Range("M14").Value = 11.64342866666667
$ws.Range("N14").Value = 34.930286
$ws.Range("O14").Value = 0.3441590756630932
$ws.Range("P14").Value = 0.3441590756630932
$ws.Range("Q14").Value = 4789.836193053026
$ws.Range("R14").Value = 43108.52573747723
$ws.Range("S14").Value = 0.09851486681321393
$ws.Range("T14").Value = 0.09851486681321393
$ws.Range("G15").Value = 411.37678
$ws.Range("H15").Value = 1234.13034
$ws.Range("I15").Value = 0.2862480573072345
$ws.Range("J15").Value = 0.2862480573072345
$ws.Range("M15").Value = 0.1645376666666667
$ws.Range("N15").Value = 0.493613
$ws.Range("O15").Value = 0.004863441250245888
$ws.Range("P15").Value = 0.004863441250245888
$ws.Range("Q15").Value = 67.68697550204666
$ws.Range("R15").Value = 609.18277951842
$ws.Range("S15").Value = 0.001392150609710753
$ws.Range("T15").Value = 0.001392150609710753
$ws.Range("G16").Value = 411.37678
$ws.Range("H16").Value = 1234.13034
$ws.Range("I16").Value = 0.2862480573072345
$ws.Range("J16").Value = 0.2862480573072345
$ws.Range("M16").Value = 12.08065233333333
$ws.Range("N16").Value = 36.241957
$ws.Range("O16").Value = 0.3570826308533967
$ws.Range("P16").Value = 0.3570826308533967
$ws.Range("Q16").Value = 4969.699857186153
$ws.Range("R16").Value = 44727.29871467538
$ws.Range("S16").Value = 0.1022142093799411
$ws.Range("T16").Value = 0.1022142093799411
$ws.Range("G17").Value = 173.2560603333334
$ws.Range("H17").Value = 519.768181
$ws.Range("I17").Value = 0.12055666021578
$ws.Range("J17").Value = 0.12055666021578
$ws.Range("M17").Value = 5.854382333333334
$ws.Range("N17").Value = 17.563147
$ws.Range("O17").Value = 0.1730451459016118
$ws.Range("P17").Value = 0.1730451459016118
$ws.Range("Q17").Value = 1014.307218758401
$ws.Range("R17").Value = 9128.764968825608
$ws.Range("S17").Value = 0.02086174485645069
$ws.Range("T17").Value = 0.02086174485645069
$ws.Range("G18").Value = 173.2560603333334
$ws.Range("H18").Value = 519.768181
$ws.Range("I18").Value = 0.12055666021578
$ws.Range("J18").Value = 0.12055666021578
$ws.Range("O18").Value = 0.1208497063316524
$ws.Range("P18").Value = 0.1208497063316525
$ws.Range("Q18").Value = 708.3627158586836
$ws.Range("R18").Value = 6375.264442728152
$ws.Range("S18").Value = 0.01456923698340182
$ws.Range("T18").Value = 0.01456923698340182
$ws.Range("G19").Value = 173.2560603333334
$ws.Range("H19").Value = 519.768181
$ws.Range("I19").Value = 0.12055666021578
$ws.Range("J19").Value = 0.12055666021578
$ws.Range("M19").Value = 11.64342866666667
$ws.Range("N19").Value = 34.930286
$ws.Range("O19").Value = 0.3441590756630932
$ws.Range("P19").Value = 0.3441590756630932
$ws.Range("Q19").Value = 2017.294579558863
$ws.Range("R19").Value = 18155.65121602976
$ws.Range("S19").Value = 0.04149066874489243
$ws.Range("T19").Value = 0.04149066874489244
$ws.Range("G20").Value = 173.2560603333334
$ws.Range("H20").Value = 519.768181
$ws.Range("I20").Value = 0.12055666021578
$ws.Range("J20").Value = 0.12055666021578
$ws.Range("M20").Value = 0.1645376666666667
$ws.Range("N20").Value = 0.493613
$ws.Range("O20").Value = 0.004863441250245888
$ws.Range("P20").Value = 0.004863441250245888
$ws.Range("Q20").Value = 28.50714790310589
$ws.Range("R20").Value = 256.564331127953
$ws.Range("S20").Value = 0.0005863202342853016
$ws.Range("T20").Value = 0.0005863202342853018
$ws.Range("G21").Value = 173.2560603333334
$ws.Range("H21").Value = 519.768181
$ws.Range("I21").Value = 0.12055666021578
$ws.Range("J21").Value = 0.12055666021578
$ws.Range("M21").Value = 12.08065233333333
$ws.Range("N21").Value = 36.241957
$ws.Range("O21").Value = 0.3570826308533967
$ws.Range("P21").Value = 0.3570826308533967
$ws.Range("Q21").Value = 2093.046229530024
$ws.Range("R21").Value = 18837.41606577022
$ws.Range("S21").Value = 0.04304868939674974
$ws.Range("T21").Value = 0.04304868939674974
$ws.Range("G22").Value = 794.4973246666667
$ws.Range("H22").Value = 2383.491974
$ws.Range("I22").Value = 0.5528345953838922
$ws.Range("J22").Value = 0.5528345953838923
$ws.Range("M22").Value = 5.854382333333334
$ws.Range("N22").Value = 17.563147
$ws.Range("O22").Value = 0.1730451459016118
$ws.Range("P22").Value = 0.1730451459016118
$ws.Range("Q22").Value = 4651.291101409131
$ws.Range("R22").Value = 41861.61991268218
$ws.Range("S22").Value = 0.09566534321766416
$ws.Range("T22").Value = 0.09566534321766419
$ws.Range("G23").Value = 794.4973246666667
$ws.Range("H23").Value = 2383.491974
$ws.Range("I23").Value = 0.5528345953838922
$ws.Range("J23").Value = 0.5528345953838923
$ws.Range("O23").Value = 0.1208497063316524
$ws.Range("P23").Value = 0.1208497063316525
$ws.Range("Q23").Value = 3248.32667648429
$ws.Range("R23").Value = 29234.94008835861
$ws.Range("S23").Value = 0.06680989850212128
$ws.Range("T23").Value = 0.0668098985021213
$ws.Range("G24").Value = 794.4973246666667
$ws.Range("H24").Value = 2383.491974
$ws.Range("I24").Value = 0.5528345953838922
$ws.Range("J24").Value = 0.5528345953838923
$ws.Range("M24").Value = 11.64342866666667
$ws.Range("N24").Value = 34.930286
$ws.Range("O24").Value = 0.3441590756630932
$ws.Range("P24").Value = 0.3441590756630932
$ws.Range("Q24").Value = 9250.672925613841
$ws.Range("R24").Value = 83256.05633052456
$ws.Range("S24").Value = 0.1902630433419005
$ws.Range("T24").Value = 0.1902630433419005
$ws.Range("G25").Value = 794.4973246666667
$ws.Range("H25").Value = 2383.491974
$ws.Range("I25").Value = 0.5528345953838922
$ws.Range("J25").Value = 0.5528345953838923
$ws.Range("M25").Value = 0.1645376666666667
$ws.Range("N25").Value = 0.493613
$ws.Range("O25").Value = 0.004863441250245888
$ws.Range("P25").Value = 0.004863441250245888
$ws.Range("Q25").Value = 130.7247359735624
$ws.Range("R25").Value = 1176.522623762062
$ws.Range("S25").Value = 0.002688678575753016
$ws.Range("T25").Value = 0.002688678575753017
$ws.Range("G26").Value = 794.4973246666667
$ws.Range("H26").Value = 2383.491974
$ws.Range("I26").Value = 0.5528345953838922
$ws.Range("J26").Value = 0.5528345953838923
$ws.Range("M26").Value = 12.08065233333333
$ws.Range("N26").Value = 36.241957
$ws.Range("O26").Value = 0.3570826308533967
$ws.Range("P26").Value = 0.3570826308533967
$ws.Range("Q26").Value = 9598.045959061457
$ws.Range("R26").Value = 86382.41363155312
$ws.Range("S26").Value = 0.1974076317464533
$ws.Range("T26").Value = 0.1974076317464533
